# Mandatory Unique Feature row (row 36) now has a B-column score of 3.
# Dependent SUM/difference formulas in E36/G36 and the B56/D56 grand
# totals recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B36").Value = 3

# Move the active selection (no more frozen/scrolled topLeftCell).
$ws.Range("A32").Select()
